{"js": "// Highlight quantitative/impact metrics (percentages, dollar amounts, large\n// numbers) in bold + a dark slate color (2C3E50) across the resume body.\n// For each target paragraph we search for the specific metric substrings\n// and apply bold + color directly to the found sub-range; Word splits the\n// paragraph into separate runs automatically so surrounding plain text is\n// left untouched.\n\nconst HIGHLIGHT_COLOR = \"#2C3E50\";\n\nasync function highlightInParagraph(context, paragraph, needles) {\n  for (const needle of needles) {\n    const found = paragraph.search(needle, { matchCase: true, matchWholeWord: false });\n    found.load(\"items\");\n    await context.sync();\n    for (let i = 0; i < found.items.length; i++) {\n      found.items[i].font.set({ bold: true, color: HIGHLIGHT_COLOR });\n    }\n  }\n  await context.sync();\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Map of exact paragraph text -> ordered list of metric substrings to bold.\nconst targets = [\n  {\n    text: \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\",\n    needles: [\"23%\", \"64%\"],\n  },\n  {\n    text: \"\u2022 Utilized advanced sampling methods to decrease survey margin of error from \u00b14.2% to \u00b12.1%, increasing voter turnout prediction accuracy from 71% to 87%, and ensuring survey results more closely reflected true population attitudes\",\n    needles: [\"\u00b14.2%\", \"\u00b12.1%\", \"71%\", \"87%\"],\n  },\n  {\n    text: \"\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis\",\n    needles: [\"73.5%\", \"$4.7M\"],\n  },\n  {\n    text: \"\u2022 Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over $2 trillion\",\n    needles: [\"$2\"],\n  },\n  {\n    text: \"\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\",\n    needles: [\"73.5%\"],\n  },\n  {\n    text: \"\u2022 $4.7M savings enabled nonprofit access\",\n    needles: [\"$4.7M\"],\n  },\n  {\n    text: \"\u2022 178% accuracy improvement in racial classification algorithms\",\n    needles: [\"178%\"],\n  },\n];\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const match = targets.find((t) => t.text === para.text);\n  if (match) {\n    await highlightInParagraph(context, para, match.needles);\n  }\n}\n", "ps1": "# Highlight quantitative/impact metrics (percentages, dollar amounts, large\n# numbers) in bold + a dark slate color (2C3E50) across the resume body.\n# For each target paragraph we locate the specific metric substrings with\n# Find.Execute (scoped to that paragraph's Range) and apply Bold + Color to\n# the found sub-range; Word splits the paragraph into separate runs\n# automatically so surrounding plain text is left untouched.\n\n$d = $word.ActiveDocument\n\n# wdColor value for hex 2C3E50 (Word stores colors as 0x00BBGGRR).\n$HighlightColor = 5258796\n\n# Exact paragraph text (sans trailing paragraph mark) -> ordered metric\n# substrings to bold/color within that paragraph.\n$targets = @{\n    \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\" = @(\"23%\", \"64%\");\n    \"\u2022 Utilized advanced sampling methods to decrease survey margin of error from \u00b14.2% to \u00b12.1%, increasing voter turnout prediction accuracy from 71% to 87%, and ensuring survey results more closely reflected true population attitudes\" = @(\"\u00b14.2%\", \"\u00b12.1%\", \"71%\", \"87%\");\n    \"\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis\" = @(\"73.5%\", \"`$4.7M\");\n    \"\u2022 Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over `$2 trillion\" = @(\"`$2\");\n    \"\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\" = @(\"73.5%\");\n    \"\u2022 `$4.7M savings enabled nonprofit access\" = @(\"`$4.7M\");\n    \"\u2022 178% accuracy improvement in racial classification algorithms\" = @(\"178%\");\n}\n\nforeach ($p in $d.Paragraphs) {\n    $fullText = $p.Range.Text\n    $paraText = $fullText.TrimEnd([char]13, [char]7)\n    if ($targets.ContainsKey($paraText)) {\n        $needles = $targets[$paraText]\n        foreach ($needle in $needles) {\n            $r = $p.Range\n            $find = $r.Find\n            $find.ClearFormatting()\n            $find.Text = $needle\n            $find.Forward = $true\n            $find.Wrap = 0\n            $find.MatchCase = $true\n            $found = $find.Execute()\n            if ($found) {\n                $r.Font.Bold = 1\n                $r.Font.Color = $HighlightColor\n            }\n        }\n    }\n}\n"}
